# Regenerate merged AHB files
# - Rename header labels: *_old -> *_FV2210, *_new -> *_FV2304
# - Turn the data range into a native Excel Table (Table1)
# - Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newHeaders = @(
    "Segmentname_FV2210","Segmentgruppe_FV2210","Segment_FV2210","Datenelement_FV2210","Segment ID_FV2210",
    "Code_FV2210","Qualifier_FV2210","Beschreibung_FV2210","Bedingungsausdruck_FV2210","Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304","Segmentgruppe_FV2304","Segment_FV2304","Datenelement_FV2304","Segment ID_FV2304",
    "Code_FV2304","Qualifier_FV2304","Beschreibung_FV2304","Bedingungsausdruck_FV2304","Bedingung_FV2304"
)

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newHeaders[$i]
}

# Convert the used range A1:U91 into a native Excel Table named "Table1"
$range = $ws.Range("A1:U91")
$listObject = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$listObject.Name = "Table1"
$listObject.TableStyle = ""

# Freeze panes above row 2 (i.e. freeze the header row)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
